$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.129.27"
$ws.Range("E2").Value = "  -10.79%  "
$ws.Range("D3").Value = "2.298.20"
$ws.Range("E3").Value = "  -20.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "446.98"
$ws.Range("E5").Value = "  -15.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.17"
$ws.Range("E6").Value = "  -11.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.473"
$ws.Range("E8").Value = "  -14.90%  "
$ws.Range("D9").Value = "2.252.79"
$ws.Range("E9").Value = "  -22.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.38"
$ws.Range("E10").Value = "  -10.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -15.99%  "
$ws.Range("E12").Value = "  -14.93%  "
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "2.694.92"
$ws.Range("E14").Value = "  -20.85%  "
$ws.Range("D15").Value = "54.110.86"
$ws.Range("E15").Value = "  -10.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.75"
$ws.Range("E16").Value = "  -17.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000120"
$ws.Range("E17").Value = "  -14.90%  "
$ws.Range("D18").Value = "2.312.70"
$ws.Range("E18").Value = "  -20.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.05"
$ws.Range("E19").Value = "  -19.88%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.41"
$ws.Range("E20").Value = "  -19.68%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "298.62"
$ws.Range("E21").Value = "  -17.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.64"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  -20.16%  "
$ws.Range("E25").Value = "  -14.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.975"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.156"
$ws.Range("E27").Value = "  -13.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.369"
$ws.Range("E28").Value = "  -19.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("E30").Value = "  -13.90%  "
$ws.Range("D31").Value = "0.0₃0706"
$ws.Range("E31").Value = "  -18.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "145.09"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "16.87"
$ws.Range("E33").Value = "  -14.51%  "
$ws.Range("E34").Value = "  -19.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.70"
$ws.Range("E35").Value = "  -15.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.60"
$ws.Range("E36").Value = "  -18.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.842"
$ws.Range("E37").Value = "  -17.01%  "
$ws.Range("E38").Value = "  -16.64%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.10"
$ws.Range("E40").Value = "  -12.21%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("E42").Value = "  -17.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.14"
$ws.Range("E43").Value = "  -15.94%  "
$ws.Range("D44").Value = "1.925.67"
$ws.Range("E44").Value = "  -16.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0495"
$ws.Range("E45").Value = "  -15.02%  "
$ws.Range("E46").Value = "  -20.56%  "
$ws.Range("E47").Value = "  -13.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0817"
$ws.Range("E48").Value = "  -11.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.12"
$ws.Range("E49").Value = "  -21.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.05"
$ws.Range("E50").Value = "  -19.17%  "
$ws.Range("E51").Value = "  -2.73%  "
